$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.1493061553617529
$ws.Range("C2").Value = 1.037716627549764
$ws.Range("D2").Value = 2.445525955610853
$ws.Range("E2").Value = 1.563817750126546
$ws.Range("F2").Value = 1.572163572957211
$ws.Range("G2").Value = 51
$ws.Range("B3").Value = -0.0001377995884945538
$ws.Range("C3").Value = 0.8974840057777216
$ws.Range("D3").Value = 1.83869274992964
$ws.Range("E3").Value = 1.355984052240158
$ws.Range("F3").Value = 1.369750733669702
$ws.Range("G3").Value = 50
$ws.Range("B4").Value = 0.1305475558742287
$ws.Range("C4").Value = 0.9634510370349124
$ws.Range("D4").Value = 2.686943648675822
$ws.Range("E4").Value = 1.639189936729671
$ws.Range("F4").Value = 1.650916085543009
$ws.Range("G4").Value = 49
$ws.Range("B5").Value = 0.02053268214733173
$ws.Range("C5").Value = 0.8033519568783193
$ws.Range("D5").Value = 1.685360403846863
$ws.Range("E5").Value = 1.298214313527186
$ws.Range("F5").Value = 1.311788311727695
$ws.Range("G5").Value = 48
$ws.Range("B6").Value = 0.1279059253806883
$ws.Range("C6").Value = 1.069557746067552
$ws.Range("D6").Value = 2.762993934606443
$ws.Range("E6").Value = 1.662225596784757
$ws.Range("F6").Value = 1.675214447118382
$ws.Range("G6").Value = 47
$ws.Range("B7").Value = 0.06123152309397346
$ws.Range("C7").Value = 0.8594512018970224
$ws.Range("D7").Value = 1.619385894554267
$ws.Range("E7").Value = 1.272550939866168
$ws.Range("F7").Value = 1.285122418519449
$ws.Range("G7").Value = 46
$ws.Range("B8").Value = 0.1270307523583763
$ws.Range("C8").Value = 1.050670529041093
$ws.Range("D8").Value = 2.656576264181475
$ws.Range("E8").Value = 1.629900691508987
$ws.Range("F8").Value = 1.643304426965771
$ws.Range("G8").Value = 45
$ws.Range("B9").Value = 0.06820288831713031
$ws.Range("C9").Value = 0.9752875536895451
$ws.Range("D9").Value = 2.031306037983135
$ws.Range("E9").Value = 1.425238940663331
$ws.Range("F9").Value = 1.440064547781096
$ws.Range("G9").Value = 44
$ws.Range("B10").Value = 0.1569399573777606
$ws.Range("C10").Value = 1.109713135099426
$ws.Range("D10").Value = 2.53526159535855
$ws.Range("E10").Value = 1.592250481349762
$ws.Range("F10").Value = 1.603249320694302
$ws.Range("G10").Value = 43
$ws.Range("B11").Value = 0.1971152062210202
$ws.Range("C11").Value = 0.8222214692082306
$ws.Range("D11").Value = 1.593335768711567
$ws.Range("E11").Value = 1.262274046596684
$ws.Range("F11").Value = 1.261901558681323
$ws.Range("G11").Value = 42
